$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K column (column G) values, regenerated after switching the
# strikeout-count source column to "K" and recalculating std/mean and s_vals.
$kValues = @{
    2  = 2
    3  = 5
    4  = 7
    5  = 5
    6  = 3
    7  = 4
    8  = 2
    9  = 1
    10 = 2
    11 = 4
    12 = 3
    13 = 5
    14 = 1
    15 = 4
    16 = 5
    17 = 2
    18 = 5
    19 = 2
    20 = 2
    21 = 0
    22 = 0
    23 = 2
    24 = 1
    25 = 2
    26 = 0
    27 = 3
    28 = 2
    29 = 1
    30 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
